$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to force text values so Excel doesn't reinterpret numeric-looking
# strings (prices, percentages) as actual numbers. We temporarily mark the
# cell as Text ("@") so the literal string is stored, then clear the format
# back to the sheet default so no stray style survives.
function Set-Text($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 - BNB
Set-Text "D2" "305.32"
Set-Text "E2" "-4.85%"

# Row 3 - OKB
Set-Text "D3" "39.48"
Set-Text "E3" "-8.69%"

# Row 4 - HuobiToken
Set-Text "D4" "5.032"
Set-Text "E4" "-2.81%"

# Row 5 - Cronos
Set-Text "D5" "0.07669"
Set-Text "E5" "-5.98%"

# Row 6 - GateToken
Set-Text "D6" "4.251"
Set-Text "E6" "-1.74%"

# Row 7 - FTXToken
Set-Text "D7" "1.595"
Set-Text "E7" "-11.12%"

# Row 8 - MXToken
Set-Text "D8" "0.8828"
Set-Text "E8" "-7.32%"

# Row 9 - LiechtensteinCryptoassetsExchange
Set-Text "D9" "0.09698"
Set-Text "E9" "-12.65%"

# Row 10 - WazirX
Set-Text "D10" "0.1723"
Set-Text "E10" "-7.35%"

# Row 11 - BitrueCoin
Set-Text "D11" "0.04493"
Set-Text "E11" "-2.71%"

# Row 12 - MandalaExchangeToken
Set-Text "D12" "0.08894"
Set-Text "E12" "-5.38%"

# Row 13 - BitMartToken
Set-Text "D13" "0.1057"
Set-Text "E13" "-0.18%"

# Row 14 - BitForexToken
Set-Text "D14" "0.001273"
Set-Text "E14" "-1.73%"

# Row 15 - was CoinExToken, now TigerCash
Set-Text "B15" "TigerCash"
Set-Text "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-Text "D15" "0.005952"
Set-Text "E15" "-1.20%"

# Row 16 - was TigerCash, now LEO
Set-Text "B16" "LEO"
Set-Text "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-Text "D16" "3.353"
Set-Text "E16" "-0.27%"

# Row 17 - was LEO, now BTSEToken
Set-Text "B17" "BTSEToken"
Set-Text "C17" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-Text "D17" "2.462"
Set-Text "E17" "-2.08%"

# Row 18 - was BTSEToken, now BitpandaEcosystemToken
Set-Text "B18" "BitpandaEcosystemToken"
Set-Text "C18" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-Text "D18" "0.3361"
Set-Text "E18" "-0.09%"

# Row 19 - was BitpandaEcosystemToken, now MCDex
Set-Text "B19" "MCDex"
Set-Text "C19" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-Text "D19" "7.044"
Set-Text "E19" "-5.04%"

# Row 20 - was MCDex, now ProBitToken
Set-Text "B20" "ProBitToken"
Set-Text "C20" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-Text "D20" "0.1350"
Set-Text "E20" "-2.71%"

# Row 21 - was ProBitToken, now ZBToken
Set-Text "B21" "ZBToken"
Set-Text "C21" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-Text "D21" "0.3202"
Set-Text "E21" "22.05%"

# Row 22 - was ZBToken, now CoinExToken
Set-Text "B22" "CoinExToken"
Set-Text "C22" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-Text "D22" "0.04204"
Set-Text "E22" "0.54%"

# Row 23 - BitKan
Set-Text "E23" "-4.28%"

# Row 24 - HotbitToken
Set-Text "D24" "0.004064"
Set-Text "E24" "-5.73%"

# Row 25 - NitroEx
Set-Text "D25" "0.0001222"
Set-Text "E25" "9.96%"

# Row 26 - UpBots
Set-Text "E26" "0.08%"

# Row 38 - One
Set-Text "D38" "0.02327"
Set-Text "E38" "-11.78%"

# Row 39 - IDEX
Set-Text "E39" "-7.35%"

# Row 40 - KickToken
Set-Text "D40" "0.007925"
Set-Text "E40" "-0.41%"

# Row 41 - BKEXToken
Set-Text "E41" "-4.89%"

# Row 42 - Dexo
Set-Text "D42" "0.006475"
Set-Text "E42" "-1.13%"

# Row 43 - CEJI
Set-Text "D43" "0.001988"
Set-Text "E43" "-6.33%"

# Row 44 - LocalTraders
Set-Text "D44" "0.008659"
Set-Text "E44" "2.59%"

# Row 45 - PooCoin
Set-Text "D45" "0.3029"
Set-Text "E45" "-5.54%"

# Row 46 - CoinLion
Set-Text "D46" "0.00006555"
Set-Text "E46" "-6.10%"

# Row 47 - Kangarootoken
Set-Text "E47" "0.15%"

# Row 48 - CoinbaseStockToken
Set-Text "D48" "0.007011"
Set-Text "E48" "98.60%"

# Row 49 - BOLO
Set-Text "D49" "0.003372"
Set-Text "E49" "-2.68%"

# Row 50 - CryptobidCoin
Set-Text "D50" "0.00002103"
Set-Text "E50" "0.15%"

# Row 51 - SpecialPowerGold
Set-Text "D51" "0.0002003"
Set-Text "E51" "0.15%"
